$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 377.5165848720561
$ws.Range("C2").Value = 605.862092896851
$ws.Range("D2").Value = 859.4081215726277
$ws.Range("E2").Value = 1054.218203400688

$ws.Range("B3").Value = 354.2488000318697
$ws.Range("C3").Value = 576.5620703043116
$ws.Range("D3").Value = 829.5128950394774
$ws.Range("E3").Value = 1028.556352226926

$ws.Range("B4").Value = 200.3093732408394
$ws.Range("C4").Value = 347.2708387582935
$ws.Range("D4").Value = 574.0614425751155
$ws.Range("E4").Value = 770.0498650957028

$ws.Range("B5").Value = 130.3464593441165
$ws.Range("C5").Value = 219.0103028675416
$ws.Range("D5").Value = 314.0876572614094
$ws.Range("E5").Value = 372.791991144114

$ws.Range("B6").Value = 66.49382152439483
$ws.Range("C6").Value = 99.27596812055262
$ws.Range("D6").Value = 158.1941005816666
$ws.Range("E6").Value = 187.6397414329435

$ws.Range("B7").Value = 111.4180642391784
$ws.Range("C7").Value = 187.2747926392082
$ws.Range("D7").Value = 287.2662924157567
$ws.Range("E7").Value = 341.2147995783006

$ws.Range("B8").Value = 108.5610867388323
$ws.Range("C8").Value = 199.3182319901114
$ws.Range("D8").Value = 276.7555218267258
$ws.Range("E8").Value = 326.4477621688607

$ws.Range("B9").Value = 175.6678320388319
$ws.Range("C9").Value = 288.5115740790028
$ws.Range("D9").Value = 387.8404871897369
$ws.Range("E9").Value = 427.8133966681748

$ws.Range("B10").Value = 266.0637605274405
$ws.Range("C10").Value = 467.6391110285813
$ws.Range("D10").Value = 747.6112428323684
$ws.Range("E10").Value = 1009.221651457342

$ws.Range("B11").Value = 232.4492655713701
$ws.Range("C11").Value = 403.3364748608295
$ws.Range("D11").Value = 663.817445183288
$ws.Range("E11").Value = 916.8538657187961

$ws.Range("B12").Value = 56.26714088105341
$ws.Range("C12").Value = 113.2232041804238
$ws.Range("D12").Value = 283.0229133490974
$ws.Range("E12").Value = 477.0849773684319

$ws.Range("B13").Value = 203.0457098268768
$ws.Range("C13").Value = 362.3155122142903
$ws.Range("D13").Value = 623.7654678184376
$ws.Range("E13").Value = 881.0403330087587
